# #3473 swapped out two properties
# - Portfolio Manager Building ID (B6, B10) updated
# - Gross Area (SF) (L7) updated
# - active selection moved to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

$ws.Range("B6").Value = 22482006
$ws.Range("L7").Value = 134036
$ws.Range("B10").Value = 22482007

$ws.Range("D7").Select()
